$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 127
$ws1.Range("F6").Value = 245
$ws1.Range("F7").Value = 12897
$ws1.Range("F8").Value = 47
$ws1.Range("F10").Value = 230
$ws1.Range("F11").Value = 2965
$ws1.Range("F12").Value = 79
$ws1.Range("F13").Value = 6307
$ws1.Range("F16").Value = 3345
$ws1.Range("F20").Value = 27
$ws1.Range("F21").Value = 32
$ws1.Range("F23").Value = 117
$ws1.Range("F24").Value = 3571
$ws1.Range("F25").Value = 79
$ws1.Range("F27").Value = 2712
$ws1.Range("F28").Value = 392
$ws1.Range("F31").Value = 203
$ws1.Range("F32").Value = 6518
$ws1.Range("F33").Value = 15
$ws1.Range("F35").Value = 299
$ws1.Range("F36").Value = 1950
$ws1.Range("F38").Value = 86
$ws1.Range("F39").Value = 1014
$ws1.Range("F41").Value = 205
$ws1.Range("F43").Value = 1144
$ws1.Range("F44").Value = 1132
$ws1.Range("F46").Value = 1180
$ws1.Range("F47").Value = 1733
$ws1.Range("F48").Value = 149
$ws1.Range("F49").Value = 1165

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 122
$ws2.Range("F15").Value = 92
$ws2.Range("F17").Value = 12

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 420
$ws3.Range("F3").Value = 577
$ws3.Range("F4").Value = 9

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 127
$ws4.Range("F6").Value = 420
$ws4.Range("F7").Value = 577
$ws4.Range("F8").Value = 245
$ws4.Range("F10").Value = 12897
$ws4.Range("F13").Value = 230
$ws4.Range("F14").Value = 2965
$ws4.Range("F15").Value = 79
$ws4.Range("F16").Value = 6307
$ws4.Range("F18").Value = 27
$ws4.Range("F19").Value = 32
$ws4.Range("F22").Value = 122
$ws4.Range("F23").Value = 3571
$ws4.Range("F24").Value = 79
$ws4.Range("F27").Value = 2712
$ws4.Range("F30").Value = 203
$ws4.Range("F31").Value = 6518
$ws4.Range("F32").Value = 92
$ws4.Range("F34").Value = 300
$ws4.Range("F35").Value = 1950
$ws4.Range("F36").Value = 12
$ws4.Range("F38").Value = 86
$ws4.Range("F39").Value = 1014
$ws4.Range("F40").Value = 205
$ws4.Range("F42").Value = 1144
$ws4.Range("F43").Value = 1132
$ws4.Range("F45").Value = 1180
$ws4.Range("F47").Value = 1733
$ws4.Range("F48").Value = 149
$ws4.Range("F49").Value = 1165
